# Update "API Discovery & Mapping - Concerns Rating" workbook
#  - rename the raw-data sheet "Sheet1" -> "Data"
#  - turn the "Legend" sheet into a proper table with header row
#    (Column1 / Column2), and make "Legend" the active/selected sheet

$wb = $excel.ActiveWorkbook

# 1. Rename "Sheet1" to "Data"
$wsData = $wb.Worksheets.Item("Sheet1")
$wsData.Name = "Data"

# 2. Add a header row to the "Legend" sheet and turn the range into Table1
$wsLegend = $wb.Worksheets.Item("Legend")
$wsLegend.Rows.Item(1).Insert()
$wsLegend.Range("A1").Value = "Column1"
$wsLegend.Range("B1").Value = "Column2"

$tbl = $wsLegend.ListObjects.Add(1, $wsLegend.Range("A1:B7"), 0, 1)
$tbl.Name = "Table1"

# 3. Make "Legend" the active sheet/tab with the table range selected
$wsLegend.Activate()
$wsLegend.Range("A1:B7").Select() | Out-Null
